$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "43.159.06"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "2.380.44"
$ws.Range("E3").Value = "  +7.01%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.50%  "
$ws.Range("E7").Value = "  +2.53%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.654"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.59%  "
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "2.734.13"
$ws.Range("E16").Value = "  +6.76%  "
$ws.Range("D17").Value = "2.384.44"
$ws.Range("E17").Value = "  +7.94%  "
$ws.Range("D18").Value = "43.121.31"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.49%  "
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "277.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.53%  "
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.17%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.13%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0919"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("E35").Value = "  +5.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0364"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("E40").Value = "  +18.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +21.60%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.230"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +21.24%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +65.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.92%  "
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.596.11"
$ws.Range("E51").Value = "  +11.26%  "
